# failurerates_v1.xlsx update
# - finalize distribution for DC system with passive balancing (X4, P4)
#   includes PB systems, DC/DC converters, and inverters
# - add two new shared strings: "Passive balancing", "system"
# - add a new block of rows (17-22, 24) mirroring the existing A:F table,
#   but using Count=1 for diode/mosfet, a new resistor row, and a combined
#   "system" total row computed with SUM(B19:B22*C19:C22)
# - fix M4 shared-formula typo (cosmetic) and move selection to E27

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section header -------------------------------------------------------
$ws.Range("A17").Value = "Passive balancing"

# --- Column headers (row 18), same headers as row 1 -----------------------
$ws.Range("A18").Value = "Component"
$ws.Range("B18").Value = "Count"
$ws.Range("C18").Value = [char]0x03BB + " (failures/10^6 hours)"
$ws.Range("D18").Value = [char]0x03BB + "(fails/hour)"
$ws.Range("E18").Value = "system lifetime"
$ws.Range("F18").Value = "R(system life)"

# --- Row 19: Capacitor -----------------------------------------------------
$ws.Range("A19").Value = "Capacitor"
$ws.Range("B19").Value = 1
$ws.Range("C19").Value = 0.001579
$ws.Range("E19").Value = 96000
$ws.Range("D19").Formula = "=C19/10^6"
$ws.Range("F19").Formula = "=EXP(-D19*E19)"
$ws.Range("F19").NumberFormat = "0.000000000"

# --- Row 20: diode -----------------------------------------------------
$ws.Range("A20").Value = "diode"
$ws.Range("B20").Value = 1
$ws.Range("C20").Value = 0.015662
$ws.Range("E20").Value = 96000

# --- Row 21: mosfet -----------------------------------------------------
$ws.Range("A21").Value = "mosfet"
$ws.Range("B21").Value = 1
$ws.Range("C21").Value = 0.075132
$ws.Range("E21").Value = 96000

# --- Row 22: resistor -----------------------------------------------------
$ws.Range("A22").Value = "resistor"
$ws.Range("B22").Value = 3
$ws.Range("C22").Value = 0.00264
$ws.Range("E22").Value = 96000

# --- Row 24: combined "system" total (note: row 23 is intentionally left
#     blank/absent, matching the source table) -----------------------------
$ws.Range("B24").Value = "system"
$ws.Range("E24").Value = 96000

# Fill D20:D24 (lambda/hour) as one shared formula group, then remove the
# accidental D23 cell so the row stays absent, just like the target sheet.
$ws.Range("D20:D24").Formula = "=C20/10^6"
$ws.Range("D23").ClearContents()

# Fill F20:F22 (reliability) as a shared formula group.
$ws.Range("F20:F22").Formula = "=EXP(-D20*E20)"
$ws.Range("F20:F22").NumberFormat = "0.000000000"

# C24 combines like terms via an array-entered SUM of element-wise products.
$ws.Range("C24").FormulaArray = "=SUM(B19:B22*C19:C22)"

# F24 (its own shared-formula group since it is not contiguous with F20:F22)
$ws.Range("F24").Formula = "=EXP(-D24*E24)"
$ws.Range("F24").NumberFormat = "0.000000000"

# --- Move the active selection to E27, like in the saved workbook ---------
$ws.Range("E27").Select() | Out-Null

# --- Explicit portrait page orientation (adds <pageSetup .../>) -----------
$ws.PageSetup.Orientation = 1
